$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.132.64'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '3.271.35'
$ws.Range("E3").Value = '  +0.78%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.01'
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.46'
$ws.Range("E6").Value = '  +2.06%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.603'
$ws.Range("E8").Value = '  +1.45%  '
$ws.Range("E9").Value = '  -2.27%  '
$ws.Range("E10").Value = '  -0.52%  '
$ws.Range("E11").Value = '  -2.24%  '
$ws.Range("D12").Value = '3.840.71'
$ws.Range("E12").Value = '  +0.76%  '
$ws.Range("E13").Value = '  +0.91%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.53'
$ws.Range("E14").Value = '  -2.09%  '
$ws.Range("D15").Value = '68.140.88'
$ws.Range("E15").Value = '  -0.08%  '
$ws.Range("E16").Value = '  -1.26%  '
$ws.Range("D17").Value = '3.345.09'
$ws.Range("E17").Value = '  +3.24%  '
$ws.Range("E18").Value = '  -0.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.35'
$ws.Range("E19").Value = '  -0.76%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '417.62'
$ws.Range("E20").Value = '  +6.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.56'
$ws.Range("E21").Value = '  -0.93%  '
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.41'
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.508'
$ws.Range("E24").Value = '  -1.13%  '
$ws.Range("E25").Value = '  -0.76%  '
$ws.Range("E26").Value = '  -0.99%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.44'
$ws.Range("E27").Value = '  -1.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.73'
$ws.Range("E30").Value = '  -0.98%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.47'
$ws.Range("E31").Value = '  -3.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.90'
$ws.Range("E32").Value = '  -2.80%  '
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("E34").Value = '  -1.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '163.89'
$ws.Range("E35").Value = '  -0.16%  '
$ws.Range("E36").Value = '  -2.20%  '
$ws.Range("E37").Value = '  -1.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '27.20'
$ws.Range("E38").Value = '  +3.79%  '
$ws.Range("E39").Value = '  -2.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.33'
$ws.Range("E41").Value = '  -3.48%  '
$ws.Range("D42").Value = '2.665.01'
$ws.Range("E42").Value = '  +3.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.91'
$ws.Range("E43").Value = '  -0.90%  '
$ws.Range("E44").Value = '  -1.08%  '
$ws.Range("E45").Value = '  -1.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '338.03'
$ws.Range("E46").Value = '  -0.87%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.52'
$ws.Range("E47").Value = '  -0.26%  '
$ws.Range("E48").Value = '  -2.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.28'
$ws.Range("E49").Value = '  -0.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.979'
$ws.Range("E50").Value = '  +0.26%  '
$ws.Range("E51").Value = '  -1.21%  '
